# Update the "Förändrad" date column (C) for all data rows from 45204 to 45205
# (2023-10-05 -> 2023-10-06) across the whole data range (rows 2 through 490).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 490 }

$range = $ws.Range("C2:C$lastRow")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
